$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("full_signals - with decay")

# 1. Insert a new column before column D (shifts D:O -> E:P, formulas/merges auto-adjust)
$ws.Range("D1").EntireColumn.Insert()

# 2. New column D content: header + "70/30" / "80/20" split values
$ws.Range("D5").Value = "train y val (%)"
$ws.Range("D6").Value = "70/30"
$ws.Range("D7").Value = "70/30"
$ws.Range("D8").Value = "70/30"
$ws.Range("D9").Value = "80/20"

# 3. New row 9 data: "modelo 4 - norm z-core"
$ws.Range("B9").Value = "modelo 4 - norm z-core"
$ws.Range("C9").Value = 30
$ws.Range("E9").Value = 0.00001
$ws.Range("F9").Value = 500
$ws.Range("G9").Value = 8
$ws.Range("H9").Formula = "=80%*F9"
$ws.Range("I9").Value = "Adam"

# 4. Fill in row 8 metrics that were previously blank
$ws.Range("K8").Value = 6.4745
$ws.Range("L8").Value = 0.3661
$ws.Range("M8").Value = 5.7979
$ws.Range("N8").Value = 0.3093
$ws.Range("P8").Value = 116

# 5. Column width tweaks (best effort given engine rounding)
$ws.Columns.Item(2).ColumnWidth = 21.736979166666668
$ws.Columns.Item(9).ColumnWidth = 8.166666666666666
$ws.Columns.Item(14).ColumnWidth = 10.307291666666666

# 6. Update selection to match author's final cursor position
$ws.Range("K3").Select()
